$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.212.57"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "2.334.32"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "2.330.79"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "2.746.72"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "60.144.36"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "2.331.54"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.33%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  +6.98%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  +11.96%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("E34").Value = "  +12.39%  "
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "320.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0496"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.560"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "0.0₆0212"
$ws.Range("E51").Value = "  +14.12%  "
